$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.408.59"
$ws.Range("E2").Value = "  +0.52%  "
$ws.Range("D3").Value = "1.577.64"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.60"
$ws.Range("E5").Value = "  +0.96%  "
$ws.Range("E6").Value = "  +0.35%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.65"
$ws.Range("E8").Value = "  -3.36%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "23.81"
$ws.Range("E9").Value = "  +0.57%  "
$ws.Range("E10").Value = "  -0.13%  "
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0895"
$ws.Range("E12").Value = "  +1.75%  "
$ws.Range("D13").Value = "1.802.65"
$ws.Range("E13").Value = "  +0.03%  "
$ws.Range("D14").Value = "1.589.66"
$ws.Range("E14").Value = "  +1.06%  "
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("D16").Value = "28.410.06"
$ws.Range("E16").Value = "  +0.37%  "
$ws.Range("E17").Value = "  -0.97%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.79"
$ws.Range("E18").Value = "  -0.85%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "229.98"
$ws.Range("E19").Value = "  +1.35%  "
$ws.Range("E20").Value = "  +0.87%  "
$ws.Range("E21").Value = "  -0.78%  "
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.94"
$ws.Range("E23").Value = "  +0.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.08"
$ws.Range("E24").Value = "  -0.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.03"
$ws.Range("E25").Value = "  +1.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.79"
$ws.Range("E26").Value = "  +0.78%  "
$ws.Range("E27").Value = "  +0.33%  "
$ws.Range("E28").Value = "  -0.64%  "
$ws.Range("E29").Value = "  -0.34%  "
$ws.Range("E30").Value = "  -0.13%  "
$ws.Range("E31").Value = "  +4.24%  "
$ws.Range("E32").Value = "  -2.07%  "
$ws.Range("E33").Value = "  +0.62%  "
$ws.Range("E34").Value = "  -0.43%  "
$ws.Range("D35").Value = "1.395.41"
$ws.Range("E35").Value = "  +0.92%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.08"
$ws.Range("E36").Value = "  +7.61%  "
$ws.Range("E37").Value = "  -3.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.36"
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("E39").Value = "  +2.43%  "
$ws.Range("E40").Value = "  -0.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.520"
$ws.Range("E41").Value = "  -2.10%  "
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("E43").Value = "  +2.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.786"
$ws.Range("E44").Value = "  -0.59%  "
$ws.Range("E45").Value = "  -3.25%  "
$ws.Range("E46").Value = "  -2.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.929"
$ws.Range("E47").Value = "  -5.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "62.56"
$ws.Range("E48").Value = "  +1.32%  "
$ws.Range("D49").Value = "1.715.25"
$ws.Range("E49").Value = "  +0.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "85.88"
$ws.Range("E50").Value = "  +0.15%  "
$ws.Range("E51").Value = "  -1.77%  "
